# This workbook tracks daily/weekly price reports for "Haba" (Hortaliza)
# at "Feria Lagunitas de Puerto Montt". A new, more recent report was
# added to the top of the data block (just under the fixed first 69
# rows), pushing the existing reports (originally rows 70-134) down by
# one row (to rows 71-135).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 70; this shifts rows 70-134 down to
# rows 71-135, carrying their values/styles/number formats with them
# (matching the diff, where every existing record below row 69 moves
# down by exactly one row).
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 4
$ws.Range("B70").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C70").Value = "Los Lagos"
$ws.Range("D70").Value = 45090
$ws.Range("E70").Value = 10
$ws.Range("F70").Value = 100112026
$ws.Range("G70").Value = "Haba"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 100
$ws.Range("K70").Value = 19000
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = 19500
$ws.Range("N70").Value = "$/saco 25 kilos"
$ws.Range("O70").Value = "Provincia de Limarí"
$ws.Range("P70").Value = 780
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
